$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.552.38"
$ws.Range("E2").Value = "  -1.43%  "
$ws.Range("D3").Value = "2.288.33"
$ws.Range("E3").Value = "  +1.39%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "95.64"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.28%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "267.54"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.62%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.622"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.87%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.604"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -5.92%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "44.90"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -6.40%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0936"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.79%  "
$ws.Range("E12").Value = "  -5.02%  "
$ws.Range("E13").Value = "  +0.21%  "
$ws.Range("D14").Value = "2.631.11"
$ws.Range("E14").Value = "  +1.37%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.16"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.75%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.851"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.68%  "
$ws.Range("D17").Value = "2.286.08"
$ws.Range("E17").Value = "  +1.92%  "
$ws.Range("D18").Value = "43.577.95"
$ws.Range("E18").Value = "  -1.31%  "
$ws.Range("E19").Value = "  +1.16%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.17"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.72%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "71.93"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.44%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.44"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.56%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "234.94"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("E24").Value = "  -12.78%  "
$ws.Range("E25").Value = "  +0.02%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.49"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.29%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.18"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.24%  "
$ws.Range("E28").Value = "  +3.98%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "40.47"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.49%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.28"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.05%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "175.16"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.21%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.92"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.80%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0882"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.29%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.35"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -6.58%  "
$ws.Range("E35").Value = "  +0.26%  "
$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.108"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.28%  "
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0355"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.37%  "
$ws.Range("E38").Value = "  -1.12%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.33"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -7.90%  "
$ws.Range("E40").Value = "  -6.98%  "
$ws.Range("E41").Value = "  +6.86%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "12.15"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.69%  "
$ws.Range("E43").Value = "  +15.56%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "63.32"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.10%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.80"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.82%  "
$ws.Range("E46").Value = "  -4.85%  "
$ws.Range("E47").Value = "  -1.30%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "98.09"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.27%  "
$ws.Range("E49").Value = "  +0.06%  "
$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").Value = "2.509.66"
$ws.Range("E50").Value = "  +1.35%  "
$ws.Range("B51").Value = "WOONetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.426"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.38%  "
